$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-01 Tuesday", "2024-10-02 Wednesday"),
    @("492÷3=", "165÷4="),
    @("208÷3=", "187÷3="),
    @("993÷3=", "971÷3="),
    @("944÷7=", "235÷8="),
    @("808÷2=", "726÷8="),
    @("689÷9=", "948÷6="),
    @("756÷3=", "141÷5="),
    @("532÷5=", "704÷9="),
    @("442÷2=", "186÷7="),
    @("683÷6=", "419÷7="),
    @("670÷8=", "270÷7="),
    @("360÷7=", "239÷9="),
    @("599÷7=", "843÷2="),
    @("688÷2=", "579÷2="),
    @("284÷8=", "100÷2="),
    @("404÷2=", "457÷4="),
    @("259÷8=", "167÷3="),
    @("978÷8=", "508÷4="),
    @("731÷3=", "867÷7="),
    @("168÷4=", "167÷9="),
    @("877÷2=", "770÷2="),
    @("292÷2=", "214÷7="),
    @("236÷6=", "744÷3="),
    @("720÷5=", "681÷4="),
    @("445÷5=", "577÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
